$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data from the GitHub Actions refresh
# Cells whose new value is a plain decimal number (e.g. "323.99") need to be
# forced to Text so Excel does not silently convert them to numeric values
# (which would also drop meaningful trailing zeros like "20.00" -> 20).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$ws.Range('D2').Value = '47.319.35'
$ws.Range('E2').Value = '  +2.93%  '
$ws.Range('D3').Value = '2.503.12'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('E4').Value = '  +0.17%  '
Set-TextValue 'D5' '323.99'
$ws.Range('E5').Value = '  +0.65%  '
Set-TextValue 'D6' '109.21'
$ws.Range('E6').Value = '  +4.24%  '
$ws.Range('E7').Value = '  +1.55%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  +0.59%  '
Set-TextValue 'D10' '38.91'
$ws.Range('E10').Value = '  +7.90%  '
$ws.Range('E11').Value = '  +1.26%  '
Set-TextValue 'D13' '18.52'
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('D15').Value = '2.893.06'
$ws.Range('E15').Value = '  +2.59%  '
$ws.Range('D16').Value = '2.502.91'
$ws.Range('E16').Value = '  +2.89%  '
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('D18').Value = '47.273.61'
$ws.Range('E18').Value = '  +3.32%  '
$ws.Range('E19').Value = '  +4.42%  '
Set-TextValue 'D20' '6.75'
$ws.Range('E20').Value = '  +4.89%  '
$ws.Range('D21').Value = '0.0₃0946'
$ws.Range('E21').Value = '  +1.70%  '
Set-TextValue 'D22' '71.08'
$ws.Range('E22').Value = '  -0.79%  '
Set-TextValue 'D23' '2.57'
$ws.Range('E23').Value = '  +8.45%  '
Set-TextValue 'D24' '250.08'
$ws.Range('E24').Value = '  +1.12%  '
Set-TextValue 'D25' '2.59'
$ws.Range('E25').Value = '  +3.54%  '
Set-TextValue 'D26' '26.16'
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +4.93%  '
Set-TextValue 'D29' '10.02'
$ws.Range('E29').Value = '  +3.56%  '
Set-TextValue 'D30' '35.77'
$ws.Range('E30').Value = '  +6.61%  '
$ws.Range('E31').Value = '  +5.80%  '
Set-TextValue 'D32' '49.93'
$ws.Range('E32').Value = '  +1.13%  '
Set-TextValue 'D33' '20.00'
$ws.Range('E33').Value = '  -1.73%  '
$ws.Range('E34').Value = '  +3.52%  '
$ws.Range('E35').Value = '  +4.44%  '
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  +5.09%  '
Set-TextValue 'D38' '4.73'
$ws.Range('E38').Value = '  +3.91%  '
Set-TextValue 'D39' '3.00'
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D41' '2.25'
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D42' '121.96'
$ws.Range('E42').Value = '  -3.87%  '
Set-TextValue 'D43' '21.32'
$ws.Range('E43').Value = '  +2.26%  '
$ws.Range('E44').Value = '  +2.46%  '
$ws.Range('D45').Value = '1.988.59'
$ws.Range('E45').Value = '  +1.32%  '
Set-TextValue 'D46' '3.06'
$ws.Range('E46').Value = '  +3.12%  '
$ws.Range('E48').Value = '  -2.47%  '
Set-TextValue 'D49' '9.07'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('E50').Value = '  +9.56%  '
Set-TextValue 'D51' '78.37'
$ws.Range('E51').Value = '  +0.78%  '
